$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 5 and 6 were previously a blank template row (old row 6, styled
# A:K). Bring that formatting down into the new row 5 (A:J - this row
# never carries a value in column K) before filling in the real data,
# then fill in the template's own row (row 6) with its data too.
$ws.Range("A6:J6").Copy()
$ws.Range("A5:J5").PasteSpecial(-4122)
$ws.Rows.Item(5).RowHeight = 18

# Write column-by-column (not row-by-row) so new shared-string entries
# land in the same order the source workbook used: both GEF names, then
# both exporter names.
$ws.Range("A5").Value = "Crumpet GEF"
$ws.Range("A6").Value = "Scone GEF"

$ws.Range("B5").Value = 20001371
$ws.Range("B6").Value = 20001371

$ws.Range("C5").Value = "Crumpet exporter"
$ws.Range("C6").Value = "Scone exporter"

$ws.Range("D5").Value = "GBP"
$ws.Range("D6").Value = "GBP"

$ws.Range("E5").Value = 7000000
$ws.Range("E6").Value = 770000

$ws.Range("F5").Value = 3938753.8
$ws.Range("F6").Value = 761579.37

$ws.Range("G5").Value = 777
$ws.Range("G6").Value = 777

$ws.Range("H5").Value = 456
$ws.Range("H6").Value = 456.77

$ws.Range("I5").Value = "GBP"
$ws.Range("I6").Value = "GBP"

$ws.Range("J5").Value = "GBP"
$ws.Range("J6").Value = "GBP"

# Move the active cell / selection to F8 (was G2).
[void]$ws.Range("F8").Select()
